$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 16.04841621693257
$ws.Range("D2").Value = 10.85818944006654
$ws.Range("E2").Value = 17.09177894305614
$ws.Range("F2").Value = 32.76323992991206
$ws.Range("G2").Value = 32.62068871839873
$ws.Range("H2").Value = 15.38480504189943
$ws.Range("J2").Value = 12.045331198532
$ws.Range("K2").Value = 9.519875163771692
$ws.Range("L2").Value = 8.204982860043238
$ws.Range("O2").Value = 23.90747917842256
# Row 3
$ws.Range("B3").Value = 15.92220147595098
$ws.Range("D3").Value = 10.86907928926056
$ws.Range("E3").Value = 17.14693553318524
$ws.Range("F3").Value = 32.87121104216952
$ws.Range("G3").Value = 32.74230945362521
$ws.Range("H3").Value = 15.43622466408308
$ws.Range("J3").Value = 12.07730703476487
$ws.Range("K3").Value = 9.199464448298347
$ws.Range("L3").Value = 8.141740518443177
$ws.Range("O3").Value = 23.99704355473635
# Row 4
$ws.Range("B4").Value = 15.84669702643472
$ws.Range("D4").Value = 10.87726777809567
$ws.Range("E4").Value = 17.18298980587958
$ws.Range("F4").Value = 32.94438661728517
$ws.Range("G4").Value = 32.82609157343225
$ws.Range("H4").Value = 15.47003743215554
$ws.Range("J4").Value = 12.09798229762638
$ws.Range("K4").Value = 8.995565415516671
$ws.Range("L4").Value = 8.103301633498768
$ws.Range("O4").Value = 24.05663077746995
# Row 5
$ws.Range("B5").Value = 15.81645656211388
$ws.Range("D5").Value = 10.88098297323938
$ws.Range("E5").Value = 17.19823330068785
$ws.Range("F5").Value = 32.97593412914917
$ws.Range("G5").Value = 32.86251575330279
$ws.Range("H5").Value = 15.48438030294577
$ws.Range("J5").Value = 12.10667039176681
$ws.Range("K5").Value = 8.910757736182532
$ws.Range("L5").Value = 8.087747532067262
$ws.Range("O5").Value = 24.08206741904582
# Row 6
$ws.Range("B6").Value = 15.81146784014791
$ws.Range("D6").Value = 10.88162274467964
$ws.Range("E6").Value = 17.20079778279391
$ws.Range("F6").Value = 32.98127685659804
$ws.Range("G6").Value = 32.86870157578161
$ws.Range("H6").Value = 15.48679599532353
$ws.Range("J6").Value = 12.10812893496122
$ws.Range("K6").Value = 8.896574421335227
$ws.Range("L6").Value = 8.085171762995532
$ws.Range("O6").Value = 24.08636084777083
# Row 7
$ws.Range("B7").Value = 15.84628701879417
$ws.Range("D7").Value = 10.87731634994981
$ws.Range("E7").Value = 17.18319315242107
$ws.Range("F7").Value = 32.94480508521718
$ws.Range("G7").Value = 32.82657357209498
$ws.Range("H7").Value = 15.47022858141285
$ws.Range("J7").Value = 12.09809840341826
$ws.Range("K7").Value = 8.994428503874568
$ws.Range("L7").Value = 8.10309140486291
$ws.Range("O7").Value = 24.05696915215907
# Row 8
$ws.Range("B8").Value = 16.00450060415811
$ws.Range("D8").Value = 10.8616328877281
$ws.Range("E8").Value = 17.11034338971604
$ws.Range("F8").Value = 32.79903885036904
$ws.Range("G8").Value = 32.66072901791976
$ws.Range("H8").Value = 15.40206966571614
$ws.Range("J8").Value = 12.05614065529813
$ws.Range("K8").Value = 9.410932843239046
$ws.Range("L8").Value = 8.183102424264153
$ws.Range("O8").Value = 23.937406876108
# Row 9
$ws.Range("B9").Value = 16.32932431563617
$ws.Range("D9").Value = 10.84276738228106
$ws.Range("E9").Value = 16.98480361428426
$ws.Range("F9").Value = 32.56789042343986
$ws.Range("G9").Value = 32.4080902432935
$ws.Range("H9").Value = 15.2861728907503
$ws.Range("J9").Value = 11.98209508902296
$ws.Range("K9").Value = 10.16756030347846
$ws.Range("L9").Value = 8.342606206243346
$ws.Range("O9").Value = 23.73943995417766
# Row 10
$ws.Range("B10").Value = 16.57515216504853
$ws.Range("D10").Value = 10.83611404811386
$ws.Range("E10").Value = 16.90306949013515
$ws.Range("F10").Value = 32.43154424841924
$ws.Range("G10").Value = 32.2671475226329
$ws.Range("H10").Value = 15.21182719254729
$ws.Range("J10").Value = 11.93266520831737
$ws.Range("K10").Value = 10.68297553056076
$ws.Range("L10").Value = 8.460691409874089
$ws.Range("O10").Value = 23.61629574451988
# Row 11
$ws.Range("B11").Value = 16.68817362710217
$ws.Range("D11").Value = 10.83464206905749
$ws.Range("E11").Value = 16.86815434679009
$ws.Range("F11").Value = 32.37681197018554
$ws.Range("G11").Value = 32.21281099919575
$ws.Range("H11").Value = 15.18034621823489
$ws.Range("J11").Value = 11.91124757625427
$ws.Range("K11").Value = 10.90800447372229
$ws.Range("L11").Value = 8.514458666528466
$ws.Range("O11").Value = 23.56512836496742
# Row 12
$ws.Range("B12").Value = 16.73111195320083
$ws.Range("D12").Value = 10.83430722605266
$ws.Range("E12").Value = 16.85525786122652
$ws.Range("F12").Value = 32.35713677457866
$ws.Range("G12").Value = 32.19364740368631
$ws.Range("H12").Value = 15.16876123733825
$ws.Range("J12").Value = 11.90329012757843
$ws.Range("K12").Value = 10.99181415608177
$ws.Range("L12").Value = 8.534813623064601
$ws.Range("O12").Value = 23.54645124527877
# Row 13
$ws.Range("B13").Value = 16.72185871967611
$ws.Range("D13").Value = 10.83436945858561
$ws.Range("E13").Value = 16.85802090199704
$ws.Range("F13").Value = 32.36132741684823
$ws.Range("G13").Value = 32.1977117128955
$ws.Range("H13").Value = 15.17124132163184
$ws.Range("J13").Value = 11.90499711494811
$ws.Range("K13").Value = 10.97382736489567
$ws.Range("L13").Value = 8.530430266356255
$ws.Range("O13").Value = 23.55044259632209
# Row 14
$ws.Range("B14").Value = 16.69170357564348
$ws.Range("D14").Value = 10.83461006679233
$ws.Range("E14").Value = 16.86708683346079
$ws.Range("F14").Value = 32.37517221037192
$ws.Range("G14").Value = 32.21120605671852
$ws.Range("H14").Value = 15.17938637844019
$ws.Range("J14").Value = 11.9105898507686
$ws.Range("K14").Value = 10.91492793079404
$ws.Range("L14").Value = 8.516133464396631
$ws.Range("O14").Value = 23.56357777411299
# Row 15
$ws.Range("B15").Value = 16.67324990239315
$ws.Range("D15").Value = 10.83478639975638
$ws.Range("E15").Value = 16.87268229679663
$ws.Range("F15").Value = 32.38378944174279
$ws.Range("G15").Value = 32.2196558537111
$ws.Range("H15").Value = 15.18441923576971
$ws.Range("J15").Value = 11.91403546090137
$ws.Range("K15").Value = 10.87866620588101
$ws.Range("L15").Value = 8.507375164175086
$ws.Range("O15").Value = 23.57171449706156
# Row 16
$ws.Range("B16").Value = 16.56778721923648
$ws.Range("D16").Value = 10.83624144588211
$ws.Range("E16").Value = 16.90539681313492
$ws.Range("F16").Value = 32.43526805530779
$ws.Range("G16").Value = 32.27089586748291
$ws.Range("H16").Value = 15.2139316026496
$ws.Range("J16").Value = 11.93408633965468
$ws.Range("K16").Value = 10.66807532574894
$ws.Range("L16").Value = 8.45717750756125
$ws.Range("O16").Value = 23.61973738253283
# Row 17
$ws.Range("B17").Value = 16.50337210377963
$ws.Range("D17").Value = 10.83753157194623
$ws.Range("E17").Value = 16.9260459838985
$ws.Range("F17").Value = 32.46871780351502
$ws.Range("G17").Value = 32.30483920983328
$ws.Range("H17").Value = 15.23263547730977
$ws.Range("J17").Value = 11.94666003345082
$ws.Range("K17").Value = 10.53643245292414
$ws.Range("L17").Value = 8.426386557068843
$ws.Range("O17").Value = 23.65044133181285
# Row 18
$ws.Range("B18").Value = 16.46643646004286
$ws.Range("D18").Value = 10.83842001439152
$ws.Range("E18").Value = 16.93813618241684
$ws.Range("F18").Value = 32.48864348304089
$ws.Range("G18").Value = 32.32528262810803
$ws.Range("H18").Value = 15.24361364465741
$ws.Range("J18").Value = 11.95399267983905
$ws.Range("K18").Value = 10.45982919339709
$ws.Range("L18").Value = 8.408681851272378
$ws.Range("O18").Value = 23.66855801011297
# Row 19
$ws.Range("B19").Value = 16.45395127929763
$ws.Range("D19").Value = 10.83874599839452
$ws.Range("E19").Value = 16.94226638038518
$ws.Range("F19").Value = 32.49550777590739
$ws.Range("G19").Value = 32.33236224845196
$ws.Range("H19").Value = 15.24736848820686
$ws.Range("J19").Value = 11.95649268496807
$ws.Range("K19").Value = 10.4337420481279
$ws.Range("L19").Value = 8.402688663038576
$ws.Range("O19").Value = 23.67477038667589
# Row 20
$ws.Range("B20").Value = 16.51021761258971
$ws.Range("D20").Value = 10.83737909134003
$ws.Range("E20").Value = 16.92382576979103
$ws.Range("F20").Value = 32.46508597858192
$ws.Range("G20").Value = 32.30113060869941
$ws.Range("H20").Value = 15.23062163022201
$ws.Range("J20").Value = 11.94531113645282
$ws.Range("K20").Value = 10.55053809409445
$ws.Range("L20").Value = 8.429663835963092
$ws.Range("O20").Value = 23.64712558095474
# Row 21
$ws.Range("B21").Value = 16.70055734024563
$ws.Range("D21").Value = 10.83453336233494
$ws.Range("E21").Value = 16.86441513179484
$ws.Range("F21").Value = 32.37107712649361
$ws.Range("G21").Value = 32.20720405958373
$ws.Range("H21").Value = 15.17698485552412
$ws.Range("J21").Value = 11.90894298316301
$ws.Range("K21").Value = 10.9322665685015
$ws.Range("L21").Value = 8.520333032661201
$ws.Range("O21").Value = 23.55970067769695
# Row 22
$ws.Range("B22").Value = 16.82575394119933
$ws.Range("D22").Value = 10.83397026738196
$ws.Range("E22").Value = 16.82748156541849
$ws.Range("F22").Value = 32.31576221477
$ws.Range("G22").Value = 32.15405310149831
$ws.Range("H22").Value = 15.14388955495085
$ws.Range("J22").Value = 11.8860654584927
$ws.Range("K22").Value = 11.17354753653752
$ws.Range("L22").Value = 8.579552501721114
$ws.Range("O22").Value = 23.50663746449424
# Row 23
$ws.Range("B23").Value = 16.75887166137651
$ws.Range("D23").Value = 10.83415249757272
$ws.Range("E23").Value = 16.84702058053984
$ws.Range("F23").Value = 32.3447237091014
$ws.Range("G23").Value = 32.18166528909137
$ws.Range("H23").Value = 15.16137391568948
$ws.Range("J23").Value = 11.89819430263484
$ws.Range("K23").Value = 11.04553570743325
$ws.Range("L23").Value = 8.547953610600306
$ws.Range("O23").Value = 23.53458513353169
# Row 24
$ws.Range("B24").Value = 16.50712245202827
$ws.Range("D24").Value = 10.8374475707063
$ws.Range("E24").Value = 16.92482884726391
$ws.Range("F24").Value = 32.46672576201004
$ws.Range("G24").Value = 32.30280437417828
$ws.Range("H24").Value = 15.23153138942254
$ws.Range("J24").Value = 11.94592064921337
$ws.Range("K24").Value = 10.5441637957711
$ws.Range("L24").Value = 8.428182185243378
$ws.Range("O24").Value = 23.64862318490595
# Row 25
$ws.Range("B25").Value = 16.24006490865622
$ws.Range("D25").Value = 10.84660208078826
$ws.Range("E25").Value = 17.01691741728855
$ws.Range("F25").Value = 32.62455202281058
$ws.Range("G25").Value = 32.46862001251983
$ws.Range("H25").Value = 15.31562704023709
$ws.Range("J25").Value = 12.00124990693696
$ws.Range("K25").Value = 9.969734687469565
$ws.Range("L25").Value = 8.299248371107511
$ws.Range("O25").Value = 23.78908263253293
